$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Preserve the "(m3/s)" shared string by writing it into the new F1
#        header cell before the old row 2 (which also referenced it) is
#        removed; this keeps the string alive & keeps its original slot.
$ws.Range("F1").Value = "(m3/s)"

# --- 2. Remove the old units/header row (old row 2). This shifts the three
#        data rows (old 3,4,5 -> Kembs, Birsfelden, Kembs-Centrale) up to
#        rows 2,3,4.
$ws.Rows.Item(2).Delete() | Out-Null

# --- 3. Fill in the rest of the new header row (row 1) left to right so the
#        shared-string table gets new entries appended in this order.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
# E1 previously held a 9pt-font header string; reset it back to the default
# (unstyled) cell format to match the rest of the new A1:E1 header cells.
$ws.Range("E1").Style = "Normal"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# --- 4. Style the whole header row with the 9pt font (matches the rest of
#        the sheet's data font) without forcing a number format override.
#        Adding a throwaway named style and deleting it afterwards leaves a
#        cellXfs entry with applyFont but not applyNumberFormat, and avoids
#        extra cellStyle/cellStyleXfs bloat.
$hdrStyle = $wb.Styles.Add("TmpHeaderStyle")
$hdrStyle.Font.Size = 9
$hdrStyle.IncludeNumber = $false
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles("TmpHeaderStyle").Delete()

# --- 5. Select A2:K2 to match the edited workbook's saved selection.
$ws.Range("A2:K2").Select() | Out-Null
